# Update version string for release "mines - version 1.0.0 (Feb 3 2026)"
$wb = $excel.ActiveWorkbook

$oldVersion = 'mines - January 30 (built on February 02 2026 12.49.33 EST)'
$newVersion = 'mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)'

# --- Sheet "About" ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Oak Grove Mine, United States, M3577, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# --- Sheet "Boundaries and methane sources" ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") rows 2-23 hold the same version string
for ($r = 2; $r -le 23; $r++) {
    $data.Cells.Item($r, 19).Value = $newVersion
}
